# EPBDS-2028 Fixed bug. Test.
# Add a new test table "testAliasTypeAsArrays(State3 state)" on the
# "Alias Datatype Usage Proper" sheet, mirroring the existing
# "testStringAliasType(State3 state)" table but using a State3[] (array)
# parameter, and update the active sheet / selections to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Build the new J3:K9 test table on sheet 2, mirroring F3:G10 ------
# Merge the header cell first so the merge operation's default border
# style doesn't clobber the header style pasted afterwards.
$ws2.Range("J3:K3").Merge()

# Copy the formatting of the existing "testStringAliasType" table
# (F3:G9 - same 7-row shape as the new table) onto the new table's
# top-left anchor so the new cells pick up matching styles (header
# style, bordered body style, etc.).
$ws2.Range("F3:G9").Copy()
$ws2.Range("J3").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("A1").Select()              # clear clipboard marquee / copy mode
$excel.CutCopyMode = $false

# Header (merged J3:K3)
$ws2.Range("J3").Value = "Rules int testAliasTypeAsArrays(State3 state)"

# Column headers
$ws2.Range("J4").Value = "C1"
$ws2.Range("K4").Value = "RET1"

# Parameter name row
$ws2.Range("J5").Value = "state"

# Parameter type row
$ws2.Range("J6").Value = "State3[]"
$ws2.Range("K6").Value = "int"

# RETURN row
$ws2.Range("J7").Value = "State"
$ws2.Range("K7").Value = "RETURN"

# Data rows
$ws2.Range("J8").Value = "CA,AR"
$ws2.Range("K8").Value = 1

$ws2.Range("J9").Value = "NY"
$ws2.Range("K9").Value = 2

# --- View / selection state -------------------------------------------
# Sheet 1 ("Alias Datatype Declaration") is no longer the active tab;
# its last selection moves to I21.
$ws1.Activate()
$ws1.Range("I21").Select()

# Sheet 2 ("Alias Datatype Usage Proper") becomes the active tab, with
# selection on I15.
$ws2.Activate()
$ws2.Range("I15").Select()
